$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value repeated for every data row.
# Find the last used row based on column A (Beteckning) and update column C
# (rows 2..lastRow) from 45204 (2023-10-05) to 45205 (2023-10-06).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$ws.Range("C2:C$lastRow").Value = 45205
